$d = $word.ActiveDocument

# 1. Remove the opening salutation paragraphs:
#      "Dear Dr. José A. Fernández Robledo,"
#      (empty paragraph)
#      "We thank you, the Editor, and the Reviewers for their constructive
#       feedback and hope this manuscript version proves acceptable for
#       publication. "
#      " "
#    leaving "Review Comments to the Author" as the new first paragraph.
$introEnd = $d.Paragraphs.Item(4).Range.End
$d.Range(0, $introEnd).Delete()

# 2. Insert a new bold paragraph "Thank you for your comment." right before
#    the "As suggested, we have more fully labelled the tables of potential
#    clade niche occupancies." paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("As suggested, we have more fully labelled the tables of potential clade niche occupancies.")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($targetIndex)
    $newPara.Range.Text = "Thank you for your comment."
    $newPara.Range.Bold = 1
}
